$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 09:52"

# Row 30: Polonia
$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 8742
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 1040
$ws.Range("E30").Value = 7355
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 347

# Row 38: Singapur
$ws.Range("A38").Value = "Singapur"
$ws.Range("B38").Value = 6588
$ws.Range("C38").Value = 596
$ws.Range("D38").Value = 740
$ws.Range("E38").Value = 5837
$ws.Range("F38").Value = 23
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 11

# Row 39: Australia
$ws.Range("A39").Value = "Australia"
$ws.Range("B39").Value = 6586
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 4167
$ws.Range("E39").Value = 2349
$ws.Range("F39").Value = 55
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 70

# Row 40: Emiratos Arabes Unidos
$ws.Range("A40").Value = "Emiratos Arabes Unidos"
$ws.Range("B40").Value = 6302
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 1188
$ws.Range("E40").Value = 5077
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 37

# Row 41: Indonesia
$ws.Range("A41").Value = "Indonesia"
$ws.Range("B41").Value = 6248
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 631
$ws.Range("E41").Value = 5082
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 535

# Row 42: Filipinas
$ws.Range("A42").Value = "Filipinas"
$ws.Range("B42").Value = 6087
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 516
$ws.Range("E42").Value = 5174
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 397

# Row 43: Serbia
$ws.Range("A43").Value = "Serbia"
$ws.Range("B43").Value = 5994
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 637
$ws.Range("E43").Value = 5240
$ws.Range("F43").Value = 126
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 117

# Row 67: Kazajistan
$ws.Range("A67").Value = "Kazajistan"
$ws.Range("B67").Value = 1661
$ws.Range("C67").Value = 46
$ws.Range("D67").Value = 382
$ws.Range("E67").Value = 1262
$ws.Range("F67").Value = 22
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 17

# Row 68: Estonia
$ws.Range("A68").Value = "Estonia"
$ws.Range("B68").Value = 1528
$ws.Range("C68").Value = 16
$ws.Range("D68").Value = 164
$ws.Range("E68").Value = 1324
$ws.Range("F68").Value = 10
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 40

# Row 69: Irak
$ws.Range("A69").Value = "Irak"
$ws.Range("B69").Value = 1513
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 953
$ws.Range("E69").Value = 478
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 82

# Row 75: Armenia
$ws.Range("A75").Value = "Armenia"
$ws.Range("B75").Value = 1291
$ws.Range("C75").Value = 43
$ws.Range("D75").Value = 545
$ws.Range("E75").Value = 726
$ws.Range("F75").Value = 30
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 20

# Row 76: Bosnia y Herzegovina
$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("B76").Value = 1268
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 338
$ws.Range("E76").Value = 883
$ws.Range("F76").Value = 4
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 47

# Row 77: Oman
$ws.Range("A77").Value = "Oman"
$ws.Range("B77").Value = 1266
$ws.Range("C77").Value = 86
$ws.Range("D77").Value = 233
$ws.Range("E77").Value = 1027
$ws.Range("F77").Value = 3
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 6
